$wb = $excel.ActiveWorkbook

# Sheet "展览" (F column "想去人数" updates)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1146
$ws1.Range("F3").Value = 584
$ws1.Range("F6").Value = 138
$ws1.Range("F10").Value = 5193
$ws1.Range("F11").Value = 4774
$ws1.Range("F14").Value = 4
$ws1.Range("F15").Value = 47

# Sheet "演出" (F column "想去人数" updates)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 75
$ws2.Range("F4").Value = 2

# Sheet "全部类型" (F column "想去人数" updates)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1146
$ws4.Range("F3").Value = 584
$ws4.Range("F6").Value = 138
$ws4.Range("F10").Value = 5193
$ws4.Range("F11").Value = 4774
$ws4.Range("F14").Value = 4
$ws4.Range("F15").Value = 47
$ws4.Range("F17").Value = 75
$ws4.Range("F19").Value = 2
